# The workbook has a "Fragen" (questions) table. The "Attribute" text for the
# "Nutzen Sie öffentliche Verkehrsmittel?" question (row 3, column D) is
# updated: the dropdown options for "Ticketart" now use square brackets
# instead of parentheses, e.g.
#   Ticketart:Dropdown(Bus,Zug,U-Bahn); ...
# becomes
#   Ticketart:Dropdown[Bus,Zug,U-Bahn]; ...
# Column D is also widened to fit the (now shorter, single-line) text, and
# the previously wrapped rows shrink back down to the default row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the "Attribute" column so the edited text fits/displays nicely.
$ws.Columns.Item(4).ColumnWidth = 88.2

# Update the cell text: swap the parentheses around the Ticketart options
# for square brackets.
$ws.Range("D3").Value = "Ticketart:Dropdown[Bus,Zug,U-Bahn]; Häufigkeit:Dropdown(Täglich,Wöchentlich,Selten)"

# Re-fit the rows whose wrapped text now only needs the default single line
# of height given the wider column.
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(5).AutoFit() | Out-Null
$ws.Rows.Item(6).AutoFit() | Out-Null

# Leave the edited cell selected, as the last user action.
$ws.Range("D3").Select() | Out-Null
